## "#5: property boat&car done"
## The 汽車 (cars) sheet (3rd sheet) gets a proper header row plus
## per-row metadata columns (property_category/category/date/
## legislator_name/legislator_id/source_file/index), matching the shape
## already used by the other property sheets (存款/股票/...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- header row (row 1): replace the stray duplicated data with real
# column labels, and extend with the same trailing metadata labels used
# on the other sheets.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# header formatting for the newly-added cells matches the existing
# bold/centered/bordered header style already used on B1:G1
$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)

# ---- data rows 2-6: add the per-row metadata columns. The existing
# A-G values (name/capacity/owner/register_date/register_reason/
# acquire_value) are left untouched.
$rows = @(2, 3, 4, 5, 6)
foreach ($r in $rows) {
    $idx = $ws.Range("A$r").Value

    $ws.Range("H$r").Value = "land"
    $ws.Range("I$r").Value = "normal"
    $ws.Range("J$r").Value = "2011-11-21"
    $ws.Range("K$r").Value = "林滄敏"
    $ws.Range("L$r").Value = 1338
    $ws.Range("M$r").Value = "tmpc7221"
    $ws.Range("N$r").Value = $idx

    $ws.Range("B$r").Copy()
    $ws.Range("H$r`:N$r").PasteSpecial(-4122)
}
